# C5-PowerPoint.pptx edit
#
# 1) The table on slide 6 switches from the deck's custom "Table_0" style
#    ({F02EB011-784E-4AD8-8680-CB5FA0212229}, defined in ppt/tableStyles.xml)
#    to the built-in gallery style {69DBBA73-0937-440F-A5D5-2842F585EC8E}
#    ("Medium Style 2 - Accent 1").
#
# 2) The theme color values carried on the deck's theme (ppt/theme/theme1.xml,
#    the theme used by the slide master / the whole deck) move from the
#    "Integral" palette to the stock "Office" palette (dk1/lt1 are already
#    identical in both palettes, so only dk2/lt2/accent1-6/hlink/folHlink
#    actually change).

$p = $ppt.ActivePresentation

# --- 1) Retarget the table's style (the table lives on slide 6, but we
#        scan every slide so this keeps working if slide order shifts) ---
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{69DBBA73-0937-440F-A5D5-2842F585EC8E}")
        }
    }
}

# --- 2) Recolor the deck theme from "Integral" to the "Office" palette ---
$scheme = $p.SlideMaster.Theme.ThemeColorScheme
$scheme.Item(3).RGB  = 6968388    # dk2      44546A
$scheme.Item(4).RGB  = 15132391   # lt2      E7E6E6
$scheme.Item(5).RGB  = 13998939   # accent1  5B9BD5
$scheme.Item(6).RGB  = 3243501    # accent2  ED7D31
$scheme.Item(7).RGB  = 10855845   # accent3  A5A5A5
$scheme.Item(8).RGB  = 49407      # accent4  FFC000
$scheme.Item(9).RGB  = 12874308   # accent5  4472C4
$scheme.Item(10).RGB = 4697456    # accent6  70AD47
$scheme.Item(11).RGB = 12673797   # hlink    0563C1
$scheme.Item(12).RGB = 7491477    # folHlink 954F72
